# Add a new "Find" / regex entry to the bash-reference sheet (工作表1),
# mirroring the existing Tar/Sed/Generic/Encode rows that live just above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 37: A=Find, B=regex, C=the long tip text (wrapped, taller row).
$findText = @'
$ find -regex ".*ABC.*"
1. option -iregex is for non capital sensitive 
2. Notice that it’s a path match but not a search (the  -name option for example is a search)
'@

$ws.Range("A37").Value = "Find"
$ws.Range("B37").Value = "regex"
$ws.Range("C37").Value = $findText

# Match the formatting used by the other long/wrapped description cells in
# column C (e.g. C31, C33-C36): wrap text + a taller row to fit the text.
$ws.Range("C37").WrapText = $true
$ws.Rows.Item(37).RowHeight = 63

# Reflect the author's final scroll/selection position after adding the row.
$ws.Range("C40").Select() | Out-Null
